$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name) from SCD0294 to SCD0018
$ws.Name = "SCD0018"

# Update the TC_ID values in column B for rows 2 and 3 (was "DGS-309")
$ws.Range("B2").Value = "SCD0018-017"
$ws.Range("B3").Value = "SCD0018-017"

# Widen column B to fit the new, longer TC_ID text (~12.57 stored character
# width; the COM layer snaps ColumnWidth to pixel increments, so feed it the
# character width that lands closest to the target after snapping)
$ws.Columns("B").ColumnWidth = 11.736979166666666

# Update the active selection to B4, as recorded in the saved view state
$ws.Range("B4").Select()
